$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3029.875
$ws.Range("I64").Value = 2927.8
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 2927.8
$ws.Range("L64").Value = 3200
$ws.Range("M64").Value = -2679.8
$ws.Range("N64").Value = -3696

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3029.875
$ws.Range("I67").Value = 2927.8
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 2927.8
$ws.Range("L67").Value = 3200
$ws.Range("M67").Value = -2069.8
$ws.Range("N67").Value = -4916

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2187.9
$ws.Range("I111").Value = 1156.3077
$ws.Range("K111").Value = 3468.9231
$ws.Range("M111").Value = -401.9231

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1220.3572
$ws.Range("J135").Value = 1843.6666
$ws.Range("L135").Value = 16592.9994
$ws.Range("N135").Value = -21662.9994

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 7947.1304
$ws.Range("J138").Value = 15453.637
$ws.Range("L138").Value = 46360.911
$ws.Range("N138").Value = -56640.911

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4672.1304
$ws.Range("I32").Value = 3528
$ws.Range("J32").Value = 7576.4614
$ws.Range("K32").Value = 3528
$ws.Range("L32").Value = 7576.4614
$ws.Range("M32").Value = -3241
$ws.Range("N32").Value = -8150.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 457426.28
$ws.Range("I122").Value = 717062.8
$ws.Range("K122").Value = 2151188.4
$ws.Range("M122").Value = -2148738.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4478.875
$ws.Range("I86").Value = 2733.3333
$ws.Range("K86").Value = 2733.3333
$ws.Range("M86").Value = -1610.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4478.875
$ws.Range("I89").Value = 2733.3333
$ws.Range("K89").Value = 13666.6665
$ws.Range("M89").Value = -8050.666499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1040.3529
$ws.Range("I105").Value = 967.875
$ws.Range("K105").Value = 967.875
$ws.Range("M105").Value = 779.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3421.8684
$ws.Range("I31").Value = 3420.0625
$ws.Range("J31").Value = 3423.182
$ws.Range("K31").Value = 3420.0625
$ws.Range("L31").Value = 3423.182
$ws.Range("M31").Value = -3125.0625
$ws.Range("N31").Value = -4013.182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3421.8684
$ws.Range("I34").Value = 3420.0625
$ws.Range("J34").Value = 3423.182
$ws.Range("K34").Value = 3420.0625
$ws.Range("L34").Value = 3423.182
$ws.Range("M34").Value = -3218.0625
$ws.Range("N34").Value = -3827.182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1041.28
$ws.Range("J107").Value = 1400
$ws.Range("L107").Value = 1400
$ws.Range("N107").Value = -5240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1690
$ws.Range("I5").Value = 785
$ws.Range("K5").Value = 2355
$ws.Range("M5").Value = -2243

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 377.77777
$ws.Range("J92").Value = 360
$ws.Range("L92").Value = 1080
$ws.Range("N92").Value = -3576

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2635.182
$ws.Range("I129").Value = 2763
$ws.Range("J129").Value = 2587.25
$ws.Range("K129").Value = 8289
$ws.Range("L129").Value = 7761.75
$ws.Range("M129").Value = -3289
$ws.Range("N129").Value = -17761.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2464.625
$ws.Range("I134").Value = 2464.625
$ws.Range("K134").Value = 7393.875
$ws.Range("M134").Value = -2323.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1690
$ws.Range("I135").Value = 785
$ws.Range("K135").Value = 7065
$ws.Range("M135").Value = -4530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 7487.8335
$ws.Range("J137").Value = 7999
$ws.Range("L137").Value = 23997
$ws.Range("N137").Value = -34197

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 4597.6
$ws.Range("I138").Value = 5124.75
$ws.Range("J138").Value = 2489
$ws.Range("K138").Value = 15374.25
$ws.Range("L138").Value = 7467
$ws.Range("M138").Value = -10234.25
$ws.Range("N138").Value = -17747

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 871.9231
$ws.Range("I140").Value = 871.9231
$ws.Range("K140").Value = 2615.7693
$ws.Range("M140").Value = 2564.2307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2751.75
$ws.Range("I80").Value = 2474.6667
$ws.Range("J80").Value = 2918
$ws.Range("K80").Value = 2474.6667
$ws.Range("L80").Value = 2918
$ws.Range("M80").Value = -1476.6667
$ws.Range("N80").Value = -4914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2751.75
$ws.Range("I83").Value = 2474.6667
$ws.Range("J83").Value = 2918
$ws.Range("K83").Value = 12373.3335
$ws.Range("L83").Value = 14590
$ws.Range("M83").Value = -7381.333500000001
$ws.Range("N83").Value = -24574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 37881.82
$ws.Range("I122").Value = 2057.6365
$ws.Range("K122").Value = 6172.9095
$ws.Range("M122").Value = -3722.9095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 31955.3
$ws.Range("J123").Value = 31955.3
$ws.Range("L123").Value = 31955.3
$ws.Range("N123").Value = -36855.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3884.2144
$ws.Range("I132").Value = 3388.111
$ws.Range("K132").Value = 10164.333
$ws.Range("M132").Value = -7634.332999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 24570.857
$ws.Range("J136").Value = 24570.857
$ws.Range("L136").Value = 73712.571
$ws.Range("N136").Value = -78812.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 32666.666
$ws.Range("J76").Value = 32666.666
$ws.Range("L76").Value = 32666.666
$ws.Range("N76").Value = -33342.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 32666.666
$ws.Range("J79").Value = 32666.666
$ws.Range("L79").Value = 32666.666
$ws.Range("N79").Value = -35006.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5256.222
$ws.Range("I132").Value = 4236.2
$ws.Range("K132").Value = 12708.6
$ws.Range("M132").Value = -10178.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3530
$ws.Range("I136").Value = 3286.5715
$ws.Range("K136").Value = 9859.7145
$ws.Range("M136").Value = -7309.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 30010
$ws.Range("J30").Value = 30010
$ws.Range("L30").Value = 30010
$ws.Range("N30").Value = -30224

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8323.375
$ws.Range("I62").Value = 7895
$ws.Range("J62").Value = 8384.571
$ws.Range("K62").Value = 7895
$ws.Range("L62").Value = 8384.571
$ws.Range("M62").Value = -7271
$ws.Range("N62").Value = -9632.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 8323.375
$ws.Range("I65").Value = 7895
$ws.Range("J65").Value = 8384.571
$ws.Range("K65").Value = 39475
$ws.Range("L65").Value = 41922.855
$ws.Range("M65").Value = -36355
$ws.Range("N65").Value = -48162.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2091.4614
$ws.Range("I126").Value = 1835.3636
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 5506.0908
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -3036.0908
$ws.Range("N126").Value = -15440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 34873.785
$ws.Range("I132").Value = 43910.137
$ws.Range("J132").Value = 1740.5
$ws.Range("K132").Value = 131730.411
$ws.Range("L132").Value = 5221.5
$ws.Range("N132").Value = -10281.5
